# Append the 2026-01-11 18:26 JST scrape run:
#  - widen column B
#  - refresh the "fetched at" timestamp on the rows that are still present
#  - insert a new top data row (row 4) for the newly scraped listing,
#    pushing the previous row 4 (cordova) down to row 5
#  - rebuild the hyperlinks so F2..F5 point at the right URLs

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Column B (title) gets wider.
$ws.Columns.Item(2).ColumnWidth = 49.1666667

# 2) Rows 2 & 3 keep their data, only the scrape timestamp changes.
$ws.Range("A2").Value = "2026-01-11 18:26:50"
$ws.Range("A3").Value = "2026-01-11 18:26:50"

# 3) Make room for the new listing: push old row 4 down to row 5.
$ws.Rows.Item(4).Insert()

# 4) New row 4 = newly scraped listing.
$ws.Range("A4").Value = "2026-01-11 18:26:50"
$ws.Range("B4").Value = "【UE5】VRoidモデルへの揺れボーン追加・PhysicsControl設定・粘液VFXの実装"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5469203"
$ws.Range("G4").Value = 18

# 5) Row 5 (previously row 4) only needs its timestamp refreshed; the rest
#    of its values/styles already moved down with the row insert above.
$ws.Range("A5").Value = "2026-01-11 18:26:50"

# 6) Rebuild hyperlinks from scratch (sheet-wide clear, then re-add in
#    order) so rIds/targets line up with the new row layout.
$ws.Range("A1:H5").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5469128")
$ws.Range("F2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5468866")
$ws.Range("F3").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5469203")
$ws.Range("F4").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5469169")
$ws.Range("F5").Style = "Hyperlink"
